$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value, as scraped from the updated coinranking.com feed.
# Values are written as literal text (matching the source sheet's inline-string
# cells for Price/Volume columns), not as numbers/percentages.
$updates = @(
    @("D2", "328.76"),
    @("E2", "0.30%"),
    @("D3", "44.22"),
    @("E3", "0.40%"),
    @("E4", "-1.06%"),
    @("D5", "0.08078"),
    @("E5", "0.47%"),
    @("D6", "2.035"),
    @("E6", "7.27%"),
    @("D7", "0.9538"),
    @("E7", "0.98%"),
    @("D8", "0.1142"),
    @("E8", "-2.25%"),
    @("D9", "0.1885"),
    @("E9", "2.40%"),
    @("E10", "2.26%"),
    @("D11", "0.09924"),
    @("E11", "2.68%"),
    @("D12", "0.04848"),
    @("E12", "11.05%"),
    @("D13", "0.1064"),
    @("E13", "-0.28%"),
    @("D14", "0.001275"),
    @("E14", "-0.10%"),
    @("D15", "0.04084"),
    @("E15", "-2.75%"),
    @("D16", "0.005824"),
    @("E16", "-2.44%"),
    @("E17", "-0.69%"),
    @("D18", "4.401"),
    @("E18", "2.87%"),
    @("E19", "2.23%"),
    @("D20", "0.3406"),
    @("E20", "-1.30%"),
    @("D21", "0.1400"),
    @("E21", "1.48%"),
    @("D22", "0.2575"),
    @("E22", "2.61%"),
    @("D23", "0.001304"),
    @("E23", "4.56%"),
    @("D24", "0.004364"),
    @("E24", "1.79%"),
    @("D25", "0.0001249"),
    @("E25", "-0.99%"),
    @("D26", "0.0003737"),
    @("E26", "-6.48%"),
    @("D38", "0.02595"),
    @("E38", "-1.94%"),
    @("D39", "0.05700"),
    @("E39", "3.48%"),
    @("D40", "0.007562"),
    @("E40", "-0.16%"),
    @("E41", "0.33%"),
    @("D42", "0.007330"),
    @("E42", "-8.76%"),
    @("D43", "0.002006"),
    @("D44", "0.009067"),
    @("E44", "2.58%"),
    @("D45", "0.00006993"),
    @("E45", "1.32%"),
    @("D46", "0.00000000749"),
    @("E46", "-0.25%"),
    @("D47", "0.0005796"),
    @("E47", "-0.26%"),
    @("D48", "0.003495"),
    @("E48", "53.71%"),
    @("D49", "0.003501"),
    @("E49", "-34.60%"),
    @("D50", "0.00002098"),
    @("E50", "-0.25%"),
    @("D51", "0.0001998"),
    @("E51", "-0.25%")
)

foreach ($update in $updates) {
    $cellRef = $update[0]
    $newValue = $update[1]
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $newValue
}
